# Update Tnc-Ptprz1 NATMI LR-pair data with refreshed TPM-based values.
# The sending/target cluster grid now covers the full 4x4 cross product of
# clusters (ECs, FAPs, MuSCs, Resolving-Mac), so the data block grows from
# 12 rows to 16 rows (plus the header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-T hold: Sending cluster, Ligand symbol, Receptor symbol,
# Target cluster, Ligand-expressing cells, Ligand detection rate,
# Ligand average/total expression value (+specificity x2), Receptor-expressing
# cells, Receptor detection rate, Receptor average/total expression value
# (+specificity x2), Edge average/total expression weight (+specificity x2).
$data = @(
    @("ECs","Tnc","Ptprz1","ECs",3,1,3.575558,10.726674,0.025194653521236,0.02519465352123599,2,0.6666666666666666,0.138766,0.416298,0.01356925767068476,0.01356925767068476,0.4961658814279999,4.465492932851999,0.0003418727455532763,0.0003418727455532762),
    @("ECs","Tnc","Ptprz1","FAPs",3,1,3.575558,10.726674,0.025194653521236,0.02519465352123599,2,0.6666666666666666,0.05416133333333333,0.162484,0.00529617548814441,0.005296175488144411,0.1936569886906666,1.742912898216,0.0001334353064114613,0.0001334353064114613),
    @("ECs","Tnc","Ptprz1","MuSCs",3,1,3.575558,10.726674,0.025194653521236,0.02519465352123599,3,1,10.002366,30.007098,0.9780831152479456,0.9780831152479456,35.764039770228,321.876357932052,0.02464246520364313,0.02464246520364312),
    @("ECs","Tnc","Ptprz1","Resolving-Mac",3,1,3.575558,10.726674,0.025194653521236,0.02519465352123599,1,0.3333333333333333,0.03120566666666667,0.09361700000000001,0.003051451593225274,0.003051451593225274,0.1115776710953333,1.004199039858,0.00007688026562813434,0.00007688026562813433),
    @("FAPs","Tnc","Ptprz1","ECs",3,1,105.8801323333333,317.640397,0.7460690747908298,0.7460690747908298,2,0.6666666666666666,0.138766,0.416298,0.01356925767068476,0.01356925767068476,14.69256244336733,132.233061990306,0.01012360351596615,0.01012360351596615),
    @("FAPs","Tnc","Ptprz1","FAPs",3,1,105.8801323333333,317.640397,0.7460690747908298,0.7460690747908298,2,0.6666666666666666,0.05416133333333333,0.162484,0.00529617548814441,0.005296175488144411,5.734609140683111,51.61148226614799,0.003951312746369772,0.003951312746369773),
    @("FAPs","Tnc","Ptprz1","MuSCs",3,1,105.8801323333333,317.640397,0.7460690747908298,0.7460690747908298,3,1,10.002366,30.007098,0.9780831152479456,0.9780831152479456,1059.051835726434,9531.466521537906,0.7297175648615674,0.7297175648615674),
    @("FAPs","Tnc","Ptprz1","Resolving-Mac",3,1,105.8801323333333,317.640397,0.7460690747908298,0.7460690747908298,1,0.3333333333333333,0.03120566666666667,0.09361700000000001,0.003051451593225274,0.003051451593225274,3.304060116216556,29.736541045949,0.002276593666926584,0.002276593666926584),
    @("MuSCs","Tnc","Ptprz1","ECs",3,1,32.36130266666667,97.08390800000001,0.2280292497513723,0.2280292497513723,2,0.6666666666666666,0.138766,0.416298,0.01356925767068476,0.01356925767068476,4.490648525842666,40.415836732584,0.003094187646329299,0.003094187646329299),
    @("MuSCs","Tnc","Ptprz1","FAPs",3,1,32.36130266666667,97.08390800000001,0.2280292497513723,0.2280292497513723,2,0.6666666666666666,0.05416133333333333,0.162484,0.00529617548814441,0.005296175488144411,1.752731300830222,15.774581707472,0.001207682923113178,0.001207682923113178),
    @("MuSCs","Tnc","Ptprz1","MuSCs",3,1,32.36130266666667,97.08390800000001,0.2280292497513723,0.2280292497513723,3,1,10.002366,30.007098,0.9780831152479456,0.9780831152479456,323.689593508776,2913.206341578984,0.223031558964474,0.223031558964474),
    @("MuSCs","Tnc","Ptprz1","Resolving-Mac",3,1,32.36130266666667,97.08390800000001,0.2280292497513723,0.2280292497513723,1,0.3333333333333333,0.03120566666666667,0.09361700000000001,0.003051451593225274,0.003051451593225274,1.009856023915111,9.088704215236001,0.0006958202174557888,0.0006958202174557888),
    @("Resolving-Mac","Tnc","Ptprz1","ECs",2,0.6666666666666666,0.1003386666666667,0.301016,0.000707021936561918,0.0007070219365619179,2,0.6666666666666666,0.138766,0.416298,0.01356925767068476,0.01356925767068476,0.01392359541866667,0.125312358768,0.000009593762836035197,0.000009593762836035197),
    @("Resolving-Mac","Tnc","Ptprz1","FAPs",2,0.6666666666666666,0.1003386666666667,0.301016,0.000707021936561918,0.0007070219365619179,2,0.6666666666666666,0.05416133333333333,0.162484,0.00529617548814441,0.005296175488144411,0.005434475971555556,0.048910283744,0.000003744512249999623,0.000003744512249999623),
    @("Resolving-Mac","Tnc","Ptprz1","MuSCs",2,0.6666666666666666,0.1003386666666667,0.301016,0.000707021936561918,0.0007070219365619179,3,1,10.002366,30.007098,0.9780831152479456,0.9780831152479456,1.003624067952,9.032616611568001,0.0006915262182611161,0.000691526218261116),
    @("Resolving-Mac","Tnc","Ptprz1","Resolving-Mac",2,0.6666666666666666,0.1003386666666667,0.301016,0.000707021936561918,0.0007070219365619179,1,0.3333333333333333,0.03120566666666667,0.09361700000000001,0.003051451593225274,0.003051451593225274,0.003131134985777778,0.028180214872,0.000002157443214767083,0.000002157443214767083)
)

# Clear the old data block (rows 2-13) before writing the refreshed,
# larger block (rows 2-17) so no stale cells are left behind.
$ws.Range("A2:T13").ClearContents()

$rowIndex = 2
foreach ($row in $data) {
    for ($col = 1; $col -le $row.Length; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $row[$col - 1]
    }
    $rowIndex++
}
